$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.904.69"
$ws.Range("E2").Value = "  +0.30%  "

$ws.Range("D3").Value = "1.883.16"
$ws.Range("E3").Value = "  +0.12%  "

$ws.Range("E4").Value = "  +1.63%  "

$ws.Range("D5").Formula = '="335.47"'
$ws.Range("D5").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4163) | Out-Null
$ws.Range("E5").Value = "  +0.59%  "

$ws.Range("D6").Formula = '="1.018"'
$ws.Range("D6").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4163) | Out-Null
$ws.Range("E6").Value = "  +1.57%  "

$ws.Range("D7").Formula = '="0.4685"'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4163) | Out-Null
$ws.Range("E7").Value = "  -1.24%  "

$ws.Range("D8").Formula = '="0.3906"'
$ws.Range("D8").Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4163) | Out-Null
$ws.Range("E8").Value = "  -1.91%  "

$ws.Range("D9").Formula = '="46.68"'
$ws.Range("D9").Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4163) | Out-Null
$ws.Range("E9").Value = "  -3.33%  "

$ws.Range("D10").Formula = '="0.07947"'
$ws.Range("D10").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4163) | Out-Null
$ws.Range("E10").Value = "  -1.24%  "

$ws.Range("E11").Value = "  -1.28%  "

$ws.Range("D12").Formula = '="21.69"'
$ws.Range("D12").Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4163) | Out-Null
$ws.Range("E12").Value = "  -1.08%  "

$ws.Range("D13").Value = "1.890.45"
$ws.Range("E13").Value = "  -0.19%  "

$ws.Range("D14").Formula = '="5.953"'
$ws.Range("D14").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4163) | Out-Null
$ws.Range("E14").Value = "  -0.27%  "

$ws.Range("E15").Value = "  -0.97%  "

$ws.Range("D16").Formula = '="1.020"'
$ws.Range("D16").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4163) | Out-Null
$ws.Range("E16").Value = "  +1.78%  "

$ws.Range("D17").Formula = '="0.06778"'
$ws.Range("D17").Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4163) | Out-Null
$ws.Range("E17").Value = "  +2.35%  "

$ws.Range("D18").Formula = '="87.40"'
$ws.Range("D18").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4163) | Out-Null
$ws.Range("E18").Value = "  +0.18%  "

$ws.Range("D19").Formula = '="0.00001047"'
$ws.Range("D19").Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4163) | Out-Null
$ws.Range("E19").Value = "  -0.63%  "

$ws.Range("D20").Formula = '="17.00"'
$ws.Range("D20").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4163) | Out-Null
$ws.Range("E20").Value = "  -2.11%  "

$ws.Range("D21").Formula = '="1.017"'
$ws.Range("D21").Copy() | Out-Null
$ws.Range("D21").PasteSpecial(-4163) | Out-Null
$ws.Range("E21").Value = "  +1.53%  "

$ws.Range("D22").Value = "27.909.78"
$ws.Range("E22").Value = "  +0.32%  "

$ws.Range("D23").Formula = '="5.465"'
$ws.Range("D23").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4163) | Out-Null
$ws.Range("E23").Value = "  -0.82%  "

$ws.Range("D24").Formula = '="10.94"'
$ws.Range("D24").Copy() | Out-Null
$ws.Range("D24").PasteSpecial(-4163) | Out-Null
$ws.Range("E24").Value = "  -1.01%  "

$ws.Range("D25").Formula = '="2.365"'
$ws.Range("D25").Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4163) | Out-Null
$ws.Range("E25").Value = "  +2.88%  "

$ws.Range("D26").Value = "2.106.63"
$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("D27").Formula = '="159.83"'
$ws.Range("D27").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4163) | Out-Null
$ws.Range("E27").Value = "  +1.73%  "

$ws.Range("D28").Formula = '="20.00"'
$ws.Range("D28").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4163) | Out-Null
$ws.Range("E28").Value = "  -1.10%  "

$ws.Range("D29").Formula = '="2.078"'
$ws.Range("D29").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4163) | Out-Null
$ws.Range("E29").Value = "  -1.41%  "

$ws.Range("E30").Value = "  -2.65%  "

$ws.Range("D31").Formula = '="120.95"'
$ws.Range("D31").Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4163) | Out-Null
$ws.Range("E31").Value = "  -1.50%  "

$ws.Range("D32").Formula = '="0.09533"'
$ws.Range("D32").Copy() | Out-Null
$ws.Range("D32").PasteSpecial(-4163) | Out-Null
$ws.Range("E32").Value = "  -0.42%  "

$ws.Range("D33").Formula = '="0.9543"'
$ws.Range("D33").Copy() | Out-Null
$ws.Range("D33").PasteSpecial(-4163) | Out-Null
$ws.Range("E33").Value = "  -2.42%  "

$ws.Range("D34").Formula = '="3.657"'
$ws.Range("D34").Copy() | Out-Null
$ws.Range("D34").PasteSpecial(-4163) | Out-Null
$ws.Range("E34").Value = "  +0.62%  "

$ws.Range("D35").Formula = '="5.327"'
$ws.Range("D35").Copy() | Out-Null
$ws.Range("D35").PasteSpecial(-4163) | Out-Null
$ws.Range("E35").Value = "  +0.21%  "

$ws.Range("E36").Value = "  -8.10%  "

$ws.Range("D37").Formula = '="0.06109"'
$ws.Range("D37").Copy() | Out-Null
$ws.Range("D37").PasteSpecial(-4163) | Out-Null
$ws.Range("E37").Value = "  -0.07%  "

$ws.Range("D38").Formula = '="0.02242"'
$ws.Range("D38").Copy() | Out-Null
$ws.Range("D38").PasteSpecial(-4163) | Out-Null
$ws.Range("E38").Value = "  -1.14%  "

$ws.Range("D39").Formula = '="1.210"'
$ws.Range("D39").Copy() | Out-Null
$ws.Range("D39").PasteSpecial(-4163) | Out-Null
$ws.Range("E39").Value = "  -1.68%  "

$ws.Range("D40").Formula = '="8.143"'
$ws.Range("D40").Copy() | Out-Null
$ws.Range("D40").PasteSpecial(-4163) | Out-Null
$ws.Range("E40").Value = "  -0.94%  "

$ws.Range("D41").Formula = '="0.5906"'
$ws.Range("D41").Copy() | Out-Null
$ws.Range("D41").PasteSpecial(-4163) | Out-Null
$ws.Range("E41").Value = "  -2.12%  "

$ws.Range("E42").Value = "  -1.02%  "

$ws.Range("E43").Value = "  -1.19%  "

$ws.Range("E44").Value = "  +1.86%  "

$ws.Range("D45").Formula = '="0.5649"'
$ws.Range("D45").Copy() | Out-Null
$ws.Range("D45").PasteSpecial(-4163) | Out-Null
$ws.Range("E45").Value = "  -1.25%  "

$ws.Range("D46").Formula = '="12.11"'
$ws.Range("D46").Copy() | Out-Null
$ws.Range("D46").PasteSpecial(-4163) | Out-Null
$ws.Range("E46").Value = "  -2.07%  "

$ws.Range("D47").Formula = '="3.400"'
$ws.Range("D47").Copy() | Out-Null
$ws.Range("D47").PasteSpecial(-4163) | Out-Null
$ws.Range("E47").Value = "  -0.49%  "

$ws.Range("D48").Formula = '="1.922"'
$ws.Range("D48").Copy() | Out-Null
$ws.Range("D48").PasteSpecial(-4163) | Out-Null
$ws.Range("E48").Value = "  -1.09%  "

$ws.Range("D49").Formula = '="0.06860"'
$ws.Range("D49").Copy() | Out-Null
$ws.Range("D49").PasteSpecial(-4163) | Out-Null
$ws.Range("E49").Value = "  +0.64%  "

$ws.Range("D50").Formula = '="113.64"'
$ws.Range("D50").Copy() | Out-Null
$ws.Range("D50").PasteSpecial(-4163) | Out-Null
$ws.Range("E50").Value = "  +0.07%  "

$ws.Range("D51").Formula = '="1.064"'
$ws.Range("D51").Copy() | Out-Null
$ws.Range("D51").PasteSpecial(-4163) | Out-Null
$ws.Range("E51").Value = "  -1.24%  "

$excel.CutCopyMode = $false
